# "Generate Report for Handoff"
#
# The two tracked e2e files swap rows on every sheet: row 2 now reports on
# 7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.md (previously row 3), and row 3 now
# reports on 08a4c95a-b82e-4aab-a37f-60064c366da2.md (previously row 2).
# The 08a4c95a... file (now row 3) also picks up a fresh handoff: its
# status becomes "Ready for handoff", its handoff/handback timestamps move
# forward, and a stale-handback warning is recorded in the Error Detail
# column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.md"
$ov.Range("B2").Value = "e2e\7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.md"

$ov.Range("A3").Value = "08a4c95a-b82e-4aab-a37f-60064c366da2.md"
$ov.Range("B3").Value = "e2e\08a4c95a-b82e-4aab-a37f-60064c366da2.md"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-23 12:48:40"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ab5096ac1f6e21fdb399b50a1b788fd43afd05f4/e2e/7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.md", "", "", "e2e\7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.md")
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ab5096ac1f6e21fdb399b50a1b788fd43afd05f4/e2e/08a4c95a-b82e-4aab-a37f-60064c366da2.md", "", "", "e2e\08a4c95a-b82e-4aab-a37f-60064c366da2.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.md"
$zh.Range("G2").Value = "7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.1b52ee52585115e8db7d0496af29e7763684c95b.zh-cn.xlf"
$zh.Range("I2").Value = "7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.md"
$zh.Range("J2").Value = "7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.1b52ee52585115e8db7d0496af29e7763684c95b.zh-cn.xlf"

$zh.Range("A3").Value = "08a4c95a-b82e-4aab-a37f-60064c366da2.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("G3").Value = "08a4c95a-b82e-4aab-a37f-60064c366da2.fc8d380b7ae84138a3ba38a5dcf1a98b08473abe.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-23 12:48:35"
$zh.Range("I3").Value = "08a4c95a-b82e-4aab-a37f-60064c366da2.md"
$zh.Range("J3").Value = "08a4c95a-b82e-4aab-a37f-60064c366da2.fc8d380b7ae84138a3ba38a5dcf1a98b08473abe.zh-cn.xlf"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ab5096ac1f6e21fdb399b50a1b788fd43afd05f4/e2e/08a4c95a-b82e-4aab-a37f-60064c366da2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4e63163807f0ec4c44e4c493db1d817839486e50/e2e/08a4c95a-b82e-4aab-a37f-60064c366da2.md."

# ColumnWidth uses Excel's character-unit scale; 39.17 round-trips to the
# OOXML raw column width of 40 (matches the target <col width="40">).
$zh.Columns.Item(16).ColumnWidth = 39.17

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ab5096ac1f6e21fdb399b50a1b788fd43afd05f4/e2e/7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.md", "", "", "7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.md")
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/2ddbd7dbf8d65961c90a3490d258c3a136f08498/e2e/7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.md", "", "", "7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.md")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ab5096ac1f6e21fdb399b50a1b788fd43afd05f4/e2e/08a4c95a-b82e-4aab-a37f-60064c366da2.md", "", "", "08a4c95a-b82e-4aab-a37f-60064c366da2.md")
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/2ddbd7dbf8d65961c90a3490d258c3a136f08498/e2e/08a4c95a-b82e-4aab-a37f-60064c366da2.md", "", "", "08a4c95a-b82e-4aab-a37f-60064c366da2.md")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.md"
$de.Range("G2").Value = "7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.1b52ee52585115e8db7d0496af29e7763684c95b.de-de.xlf"
$de.Range("I2").Value = "7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.md"
$de.Range("J2").Value = "7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.1b52ee52585115e8db7d0496af29e7763684c95b.de-de.xlf"

$de.Range("A3").Value = "08a4c95a-b82e-4aab-a37f-60064c366da2.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "08a4c95a-b82e-4aab-a37f-60064c366da2.fc8d380b7ae84138a3ba38a5dcf1a98b08473abe.de-de.xlf"
$de.Range("H3").Value = "2016-08-23 12:48:40"
$de.Range("I3").Value = "08a4c95a-b82e-4aab-a37f-60064c366da2.md"
$de.Range("J3").Value = "08a4c95a-b82e-4aab-a37f-60064c366da2.fc8d380b7ae84138a3ba38a5dcf1a98b08473abe.de-de.xlf"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ab5096ac1f6e21fdb399b50a1b788fd43afd05f4/e2e/08a4c95a-b82e-4aab-a37f-60064c366da2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4e63163807f0ec4c44e4c493db1d817839486e50/e2e/08a4c95a-b82e-4aab-a37f-60064c366da2.md."

# ColumnWidth uses Excel's character-unit scale; 39.17 round-trips to the
# OOXML raw column width of 40 (matches the target <col width="40">).
$de.Columns.Item(16).ColumnWidth = 39.17

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ab5096ac1f6e21fdb399b50a1b788fd43afd05f4/e2e/7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.md", "", "", "7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.md")
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0d37699b58a961f54095158f398c956a1826cba0/e2e/7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.md", "", "", "7c5d75a6-8973-4f2c-b30e-4e0a81cf60c5.md")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ab5096ac1f6e21fdb399b50a1b788fd43afd05f4/e2e/08a4c95a-b82e-4aab-a37f-60064c366da2.md", "", "", "08a4c95a-b82e-4aab-a37f-60064c366da2.md")
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0d37699b58a961f54095158f398c956a1826cba0/e2e/08a4c95a-b82e-4aab-a37f-60064c366da2.md", "", "", "08a4c95a-b82e-4aab-a37f-60064c366da2.md")
